$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "B"=1.02; "C"=1.045270918663713; "D"=1.053301163299191; "E"=1.05852125235453; "F"=1.065833810544642; "I"=1.043679157331845; "J"=1.050331547086078; "K"=1.056047510103612; "L"=1.061253259312229; "M"=1.068545985154074; "N"=1.020561288360648 }
  3 = @{ "B"=1.02; "C"=1.046498812723691; "D"=1.054270724515332; "E"=1.059590870449503; "F"=1.066927575217535; "I"=1.043992356127501; "J"=1.051205660242307; "K"=1.056829687670067; "L"=1.062136279030848; "M"=1.069454530318535; "N"=1.020858207392284 }
  4 = @{ "B"=1.02; "C"=1.047292946319485; "D"=1.054897543037564; "E"=1.060283401381104; "F"=1.067635300539399; "I"=1.0441931711332; "J"=1.051770364696209; "K"=1.05733463875901; "L"=1.062707468146193; "M"=1.070041810136158; "N"=1.021049861860901 }
  5 = @{ "B"=1.02; "C"=1.047626707361436; "D"=1.055160926064978; "E"=1.060574641688587; "F"=1.06793282582675; "I"=1.044277152243413; "J"=1.052007550846968; "K"=1.057546640862179; "L"=1.062947552667917; "M"=1.070288557372854; "N"=1.021130320691089 }
  6 = @{ "B"=1.02; "C"=1.047682741931958; "D"=1.055205141546468; "E"=1.060623548116765; "F"=1.067982781512098; "I"=1.044291227159497; "J"=1.052047362834444; "K"=1.05758222059066; "L"=1.062987861388747; "M"=1.07032997882694; "N"=1.021143823481816 }
  7 = @{ "B"=1.02; "C"=1.047297406412611; "D"=1.054901062890146; "E"=1.060287292553488; "F"=1.067639276095432; "I"=1.044194295027393; "J"=1.051773534834962; "K"=1.057337472638115; "L"=1.062710676337055; "M"=1.070045107754415; "N"=1.021050937398236 }
  8 = @{ "B"=1.02; "C"=1.045685974447794; "D"=1.053628945899749; "E"=1.058882648916726; "F"=1.066203456733056; "I"=1.043785386655794; "J"=1.05062714581749; "K"=1.056312093428711; "L"=1.061551718352761; "M"=1.068853158935763; "N"=1.020661731442771 }
  9 = @{ "B"=1.02; "C"=1.042843299489953; "D"=1.051383047033387; "E"=1.056410650054133; "F"=1.063673218759903; "I"=1.043050688589259; "J"=1.048600078907869; "K"=1.05449624660976; "L"=1.059508052319944; "M"=1.066748089092113; "N"=1.019972272694838 }
  10 = @{ "B"=1.02; "C"=1.040945948093946; "D"=1.049882860109249; "E"=1.054764736130864; "F"=1.061986243828227; "I"=1.042551354569129; "J"=1.047243933225962; "K"=1.053279579287949; "L"=1.058144603505371; "M"=1.065341495484183; "N"=1.0195101747408 }
  11 = @{ "B"=1.02; "C"=1.040123813526036; "D"=1.049232557055192; "E"=1.054052521748431; "F"=1.061255716331118; "I"=1.042332870258391; "J"=1.046655559316442; "K"=1.052751288673292; "L"=1.057553966905679; "M"=1.06473165036167; "N"=1.01930949338556 }
  12 = @{ "B"=1.02; "C"=1.039818347947937; "D"=1.048990897329572; "E"=1.053788044337772; "F"=1.060984356190479; "I"=1.042251373948367; "J"=1.046436836122032; "K"=1.052554836613945; "L"=1.057334539024773; "M"=1.064505008251486; "N"=1.019234862328087 }
  13 = @{ "B"=1.02; "C"=1.039883875386182; "D"=1.049042739077689; "E"=1.053844772455509; "F"=1.061042564321067; "I"=1.042268870646558; "J"=1.04648376093082; "K"=1.052596986285987; "L"=1.057381608837008; "M"=1.06455362913376; "N"=1.019250874980845 }
  14 = @{ "B"=1.02; "C"=1.040098565460983; "D"=1.049212583601532; "E"=1.054030658526099; "F"=1.061233285822442; "I"=1.042326140717918; "J"=1.046637483166529; "K"=1.052735054408354; "L"=1.057535829732747; "M"=1.064712918477866; "N"=1.01930332617857 }
  15 = @{ "B"=1.02; "C"=1.040230831333511; "D"=1.049317216041999; "E"=1.054145198411385; "F"=1.061350794310802; "I"=1.042361381428945; "J"=1.046732173282791; "K"=1.052820093361046; "L"=1.057630845100332; "M"=1.064811046165286; "N"=1.01933563127711 }
  16 = @{ "B"=1.02; "C"=1.041000498239351; "D"=1.049926003472255; "E"=1.054812013425152; "F"=1.062034725240876; "I"=1.042565806792303; "J"=1.047282957219867; "K"=1.053314609214027; "L"=1.058183796673323; "M"=1.065381952346024; "N"=1.019523480835389 }
  17 = @{ "B"=1.02; "C"=1.04148313572539; "D"=1.050307688208396; "E"=1.055230415943164; "F"=1.062463721108761; "I"=1.042693429449529; "J"=1.047628139399687; "K"=1.053624412701192; "L"=1.058530579843545; "M"=1.065739856841186; "N"=1.019641155658015 }
  18 = @{ "B"=1.02; "C"=1.041764595301238; "D"=1.050530249822687; "E"=1.055474509164085; "F"=1.062713941802709; "I"=1.042767650662279; "J"=1.047829367100005; "K"=1.053804974504363; "L"=1.058732828132239; "M"=1.065948541221986; "N"=1.019709736505565 }
  19 = @{ "B"=1.02; "C"=1.041860556546039; "D"=1.050606125976101; "E"=1.055557746494534; "F"=1.062799259725002; "I"=1.042792921061066; "J"=1.047897961726842; "K"=1.053866517477164; "L"=1.058801785450214; "M"=1.066019684473065; "N"=1.019733111166689 }
  20 = @{ "B"=1.02; "C"=1.041431358970027; "D"=1.050266744157354; "E"=1.055185520551172; "F"=1.062417694474285; "I"=1.04267975939506; "J"=1.047591116135367; "K"=1.053591188346904; "L"=1.058493375829754; "M"=1.065701464869457; "N"=1.019628536150668 }
  21 = @{ "B"=1.02; "C"=1.040035347041477; "D"=1.049162571619634; "E"=1.053975917766963; "F"=1.061177123371829; "I"=1.042309285541162; "J"=1.04659222064103; "K"=1.052694402909123; "L"=1.057490416611183; "M"=1.064666015016234; "N"=1.019287883070396 }
  22 = @{ "B"=1.02; "C"=1.039157107315078; "D"=1.048467707804029; "E"=1.053215801178474; "F"=1.060397070323719; "I"=1.042074377778367; "J"=1.04596316254029; "K"=1.052129276737408; "L"=1.056859589360621; "M"=1.064014300531744; "N"=1.019073185424504 }
  23 = @{ "B"=1.02; "C"=1.039622728219055; "D"=1.048836128073835; "E"=1.053618714951518; "F"=1.060810597005338; "I"=1.042199094332305; "J"=1.046296734675438; "K"=1.052428982697848; "L"=1.057194024655858; "M"=1.064359852135192; "N"=1.019187049723865 }
  24 = @{ "B"=1.02; "C"=1.041454754826091; "D"=1.050285245223784; "E"=1.055205806703988; "F"=1.062438491945833; "I"=1.04268593698074; "J"=1.047607845698582; "K"=1.053606201439239; "L"=1.058510186797873; "M"=1.065718812781548; "N"=1.019634238538637 }
  25 = @{ "B"=1.02; "C"=1.043578583214276; "D"=1.051964176673891; "E"=1.057049350425091; "F"=1.064327368235798; "I"=1.04324230466898; "J"=1.049124959275396; "K"=1.05496675825547; "L"=1.060036563074333; "M"=1.067292862760614; "N"=1.020150946251623 }
}

foreach ($rowKey in $data.Keys) {
  $rowData = $data[$rowKey]
  foreach ($colKey in $rowData.Keys) {
    $ws.Range("$colKey$rowKey").Value = $rowData[$colKey]
  }
}
